{"js": "// The diff fixes two genuine content/typo issues inside the \"Story points\"\n// and \"Tasks\" bullets (the surrounding w:proofErr spell/grammar-check markers\n// that disappear in the diff are just cosmetic artifacts of Word's proofing\n// engine collapsing once the text is edited, not something authored by hand).\n//\n// 1) \"...what he needs to do.\" -> \"...what he needs to do?\"\n// 2) \"...which are rwarded, when the story gets completed. Completed in the\n//     sense, that all its tasks have been completed.\" ->\n//    \"...which are rewarded, when the story gets completed. Completed in the\n//     sense, that all its tasks have been completed. These story points are\n//     given to a story, which define, how tough that story is, based on its\n//     complexity and other parameters.  We'll be using the Fibonacci Series\n//     to define these points(1,2,3,5,8,13,21.... ).. \"\n\nconst body = context.document.body;\n\n// 1) Fix the missing question mark in the \"Tasks\" bullet.\nconst doneSearch = body.search(\n  \"Which ultimately, tells the person, what he needs to do.\",\n  { matchCase: true }\n);\ndoneSearch.load(\"items\");\nawait context.sync();\n\nif (doneSearch.items.length === 0) {\n  throw new Error(\"Could not find the 'Tasks' sentence to fix.\");\n}\ndoneSearch.items[0].insertText(\n  \"Which ultimately, tells the person, what he needs to do?\",\n  Word.InsertLocation.replace\n);\n\n// 2) Fix the \"rwarded\" typo and append the missing explanation about story\n//    points / the Fibonacci series in the \"Story points\" bullet.\nconst rewardSearch = body.search(\n  \"rwarded, when the story gets completed. Completed in the sense, that all its tasks have been completed.\",\n  { matchCase: true }\n);\nrewardSearch.load(\"items\");\nawait context.sync();\n\nif (rewardSearch.items.length === 0) {\n  throw new Error(\"Could not find the 'Story points' sentence to fix.\");\n}\nrewardSearch.items[0].insertText(\n  \"rewarded, when the story gets completed. Completed in the sense, that all its tasks have been completed. These story points are given to a story, which define, how tough that story is, based on its complexity and other parameters.  We\\u2019ll be using the Fibonacci Series to define these points(1,2,3,5,8,13,21.... ).. \",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# The diff fixes two genuine content/typo issues inside the \"Story points\"\n# and \"Tasks\" bullets (the surrounding proofing/spell-check marks that\n# disappear in the diff are just cosmetic artifacts of Word's proofing\n# engine, not something authored by hand).\n#\n# 1) \"...what he needs to do.\" -> \"...what he needs to do?\"\n# 2) \"...which are rwarded, when the story gets completed. Completed in the\n#     sense, that all its tasks have been completed.\" ->\n#    \"...which are rewarded, when the story gets completed. Completed in the\n#     sense, that all its tasks have been completed. These story points are\n#     given to a story, which define, how tough that story is, based on its\n#     complexity and other parameters.  We\u2019ll be using the Fibonacci Series\n#     to define these points(1,2,3,5,8,13,21.... ).. \"\n\n$d = $word.ActiveDocument\n\n# 1) Fix the missing question mark in the \"Tasks\" bullet.\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$found1 = $find1.Execute(\n    \"Which ultimately, tells the person, what he needs to do.\",\n    $false, $true, $false, $false, $false, $true, 1, $false,\n    \"Which ultimately, tells the person, what he needs to do?\", 1\n)\nif (-not $found1) {\n    throw \"Could not find the 'Tasks' sentence to fix.\"\n}\n\n# 2) Fix the \"rwarded\" typo and append the missing explanation about story\n#    points / the Fibonacci series in the \"Story points\" bullet.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$found2 = $find2.Execute(\n    \"rwarded, when the story gets completed. Completed in the sense, that all its tasks have been completed.\",\n    $false, $true, $false, $false, $false, $true, 1, $false,\n    \"rewarded, when the story gets completed. Completed in the sense, that all its tasks have been completed. These story points are given to a story, which define, how tough that story is, based on its complexity and other parameters.  We\u2019ll be using the Fibonacci Series to define these points(1,2,3,5,8,13,21.... ).. \",\n    1\n)\nif (-not $found2) {\n    throw \"Could not find the 'Story points' sentence to fix.\"\n}\n"}
